$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (btcinr) values - quantity was fixed, recalculated dependent columns
$ws.Range("B2").Value = 3368900
$ws.Range("E2").Value = 11471.47000000001
$ws.Range("F2").Value = 50.86129459325488

# Update row 3 (ltcinr) values - quantity was fixed, recalculated dependent columns
$ws.Range("B3").Value = 13975
$ws.Range("E3").Value = 1377650
$ws.Range("F3").Value = 24.93032935215346

# Remove row 4 (ethinr) entirely - it was a bogus 0-quantity row
$ws.Range("A4:F4").EntireRow.Delete()
